$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H81").Value = 30001
$ws.Range("J81").Value = 30001
$ws.Range("L81").Value = 30001
$ws.Range("N81").Value = -31997
$ws.Range("H84").Value = 30001
$ws.Range("J84").Value = 30001
$ws.Range("L84").Value = 90003
$ws.Range("N84").Value = -99987
$ws.Range("H112").Value = 1556.25
$ws.Range("I112").Value = 1078.3334
$ws.Range("J112").Value = 2034.1666
$ws.Range("K112").Value = 3235.0002
$ws.Range("L112").Value = 6102.4998
$ws.Range("M112").Value = -2127.0002
$ws.Range("N112").Value = -8318.4998
$ws.Range("H116").Value = 113768.445
$ws.Range("I116").Value = 202682
$ws.Range("J116").Value = 2626.5
$ws.Range("K116").Value = 202682
$ws.Range("L116").Value = 2626.5
$ws.Range("M116").Value = -199240
$ws.Range("N116").Value = -9510.5
$ws.Range("H134").Value = 83333.336
$ws.Range("J134").Value = 83333.336
$ws.Range("L134").Value = 83333.336
$ws.Range("N134").Value = -93473.336
$ws.Range("H140").Value = 54850
$ws.Range("J140").Value = 54850
$ws.Range("L140").Value = 54850
$ws.Range("N140").Value = -65210

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H23").Value = 2500
$ws.Range("J23").Value = 2500
$ws.Range("L23").Value = 2500
$ws.Range("N23").Value = -3018
$ws.Range("H32").Value = 7393.45
$ws.Range("I32").Value = 7065.4126
$ws.Range("J32").Value = 18000
$ws.Range("K32").Value = 7065.4126
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = -6778.4126
$ws.Range("N32").Value = -18574
$ws.Range("H37").Value = 14750
$ws.Range("J37").Value = 14666.667
$ws.Range("L37").Value = 14666.667
$ws.Range("N37").Value = -15212.667
$ws.Range("H44").Value = 16760
$ws.Range("J44").Value = 16760
$ws.Range("L44").Value = 16760
$ws.Range("N44").Value = -17736
$ws.Range("J55").Value = 24600
$ws.Range("L55").Value = 24600
$ws.Range("N55").Value = -25230
$ws.Range("H74").Value = 4098.086
$ws.Range("I74").Value = 1099.6786
$ws.Range("K74").Value = 1099.6786
$ws.Range("M74").Value = -225.6786
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H77").Value = 4098.086
$ws.Range("I77").Value = 1099.6786
$ws.Range("K77").Value = 5498.393
$ws.Range("M77").Value = -1130.393
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1868.6666
$ws.Range("I99").Value = 1105
$ws.Range("J99").Value = 2250.5
$ws.Range("K99").Value = 1105
$ws.Range("L99").Value = 2250.5
$ws.Range("M99").Value = 393
$ws.Range("N99").Value = -5246.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 38462628
$ws.Range("I58").Value = 76923850
$ws.Range("J58").Value = 1406.4615
$ws.Range("K58").Value = 76923850
$ws.Range("L58").Value = 1406.4615
$ws.Range("M58").Value = -76923647
$ws.Range("N58").Value = -1812.4615
$ws.Range("H136").Value = 38462628
$ws.Range("I136").Value = 76923850
$ws.Range("J136").Value = 1406.4615
$ws.Range("K136").Value = 230771550
$ws.Range("L136").Value = 4219.3845
$ws.Range("M136").Value = -230769000
$ws.Range("N136").Value = -9319.3845

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3070
$ws.Range("I109").Value = 3005
$ws.Range("J109").Value = 3330
$ws.Range("K109").Value = 9015
$ws.Range("L109").Value = 9990
$ws.Range("M109").Value = -7975
$ws.Range("N109").Value = -12070
$ws.Range("H131").Value = 686.3158
$ws.Range("I131").Value = 202.5
$ws.Range("J131").Value = 1038.1818
$ws.Range("K131").Value = 607.5
$ws.Range("L131").Value = 3114.5454
$ws.Range("M131").Value = 4432.5
$ws.Range("N131").Value = -13194.5454
$ws.Range("H137").Value = 3929090.2
$ws.Range("I137").Value = 70659.375
$ws.Range("J137").Value = 10102579
$ws.Range("K137").Value = 211978.125
$ws.Range("L137").Value = 30307737
$ws.Range("M137").Value = -206878.125
$ws.Range("N137").Value = -30317937

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 52264
$ws.Range("I141").Value = 25000
$ws.Range("J141").Value = 72712
$ws.Range("K141").Value = 25000
$ws.Range("L141").Value = 72712
$ws.Range("M141").Value = -19820
$ws.Range("N141").Value = -83072

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 4425
$ws.Range("I22").Value = 2733.3333
$ws.Range("J22").Value = 9500
$ws.Range("K22").Value = 2733.3333
$ws.Range("L22").Value = 9500
$ws.Range("M22").Value = -2440.3333
$ws.Range("N22").Value = -10086
$ws.Range("H96").Value = 1510.7778
$ws.Range("I96").Value = 1199.125
$ws.Range("J96").Value = 4004
$ws.Range("K96").Value = 1199.125
$ws.Range("L96").Value = 4004
$ws.Range("M96").Value = 173.875
$ws.Range("N96").Value = -6750
$ws.Range("H135").Value = 43443.332
$ws.Range("J135").Value = 43443.332
$ws.Range("L135").Value = 43443.332
$ws.Range("N135").Value = -53583.332
$ws.Range("H140").Value = 56164.5
$ws.Range("J140").Value = 56164.5
$ws.Range("L140").Value = 56164.5
$ws.Range("N140").Value = -66524.5
$ws.Range("H141").Value = 68357.5
$ws.Range("J141").Value = 68357.5
$ws.Range("L141").Value = 68357.5
$ws.Range("N141").Value = -78717.5
